$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 13 (pushes old rows 13-23 down to 15-25),
# matching the new used range A1:C26.
$ws.Rows("13:14").Insert()

# --- Set final cell content for rows 10-26 (explicit set/clear for every cell so the
#     post-shift leftovers from the old layout are fully replaced) ---
# Row 10
$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = 'Aplicação dos conhecimentos adquiridos nas disciplinas obrigatórias e das competências desenvolvidas durante o Curso a uma situação possível do ambiente profissional'
$ws.Range("C10").Value = 'Aplicação dos conhecimentos adquiridos nas disciplinas obrigatórias e das competências desenvolvidas durante o Curso a uma situação possível do ambiente profissional'

# Row 11
$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()

# Row 12
$ws.Range("A12").Value = 'Docentes responsáveis:'
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()

# Row 13
$ws.Range("A13").ClearContents()
$ws.Range("B13").Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Range("C13").Value = '7459752 - Maria Ismenia Sodero Toledo Faria'

# Row 14
$ws.Range("A14").ClearContents()
$ws.Range("B14").Value = '2166002 - Sandra Giacomin Schneider'
$ws.Range("C14").Value = '2166002 - Sandra Giacomin Schneider'

# Row 15
$ws.Range("A15").ClearContents()
$ws.Range("B15").Value = '1922320 - Sebastiao Ribeiro'
$ws.Range("C15").Value = '1922320 - Sebastiao Ribeiro'

# Row 16
$ws.Range("A16").Value = 'Programa resumido:'
$ws.Range("B16").Value = 'Elaboração de um projeto de engenharia, de pesquisa científica ou modelo de negócio, Desenvolvimento do projeto, com características inter e transdisciplinar'
$ws.Range("C16").Value = 'Elaboração de um projeto de engenharia, de pesquisa científica ou modelo de negócio, Desenvolvimento do projeto, com características inter e transdisciplinar'

# Row 17
$ws.Range("A17").Value = 'Short syllabus:'
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()

# Row 18
$ws.Range("A18").Value = 'Programa:'
$ws.Range("B18").Value = 'A disciplina consiste no desenvolvimento de um projeto supervisionado por docente e/ou profissional de engenharia que poderá ser realizado em grupo ou de forma individual.1. Elaboração de um projeto de engenharia, ou pesquisa científica ou modelo de negócio utilizando as competências adquiridas nos Projetos de Engenharia I, II, III e IV). O projeto deve atender os princípios de planejamento e gestão de projetos ou de negócios, inclusive, se possível construindo modelo ou protótipo físico e/ou digital2. Desenvolvimento do Projeto – em projetos experimentais deverão ser produzidos alguns produtos, processos ou sistemas reais, teste de modelos ou protótipos.3. O aluno, individualmente ou em equipe, deverá elaborar uma monografia ou plano de negócio do projeto e submete-la a apreciação de uma banca – a monografia ou plano de negócio deve atender aos padrões estabelecidos e utilizados nas disciplinas de Projetos de Engenharia I, II, III e IV).4. Apresentação do Projeto Final de Curso para uma banca de três professores, sendo um orientador, no caso de ter mais de um, e dois outros membros, entre eles preferencialmente, um vindo da indústria do ramo de atividade do tema escolhido.'
$ws.Range("C18").Value = 'A disciplina consiste no desenvolvimento de um projeto supervisionado por docente e/ou profissional de engenharia que poderá ser realizado em grupo ou de forma individual.1. Elaboração de um projeto de engenharia, ou pesquisa científica ou modelo de negócio utilizando as competências adquiridas nos Projetos de Engenharia I, II, III e IV). O projeto deve atender os princípios de planejamento e gestão de projetos ou de negócios, inclusive, se possível construindo modelo ou protótipo físico e/ou digital2. Desenvolvimento do Projeto – em projetos experimentais deverão ser produzidos alguns produtos, processos ou sistemas reais, teste de modelos ou protótipos.3. O aluno, individualmente ou em equipe, deverá elaborar uma monografia ou plano de negócio do projeto e submete-la a apreciação de uma banca – a monografia ou plano de negócio deve atender aos padrões estabelecidos e utilizados nas disciplinas de Projetos de Engenharia I, II, III e IV).4. Apresentação do Projeto Final de Curso para uma banca de três professores, sendo um orientador, no caso de ter mais de um, e dois outros membros, entre eles preferencialmente, um vindo da indústria do ramo de atividade do tema escolhido.'

# Row 19
$ws.Range("A19").Value = 'Syllabus:'
$ws.Range("B19").ClearContents()
$ws.Range("C19").ClearContents()

# Row 20
$ws.Range("A20").Value = 'Avaliação:'
$ws.Range("B20").ClearContents()
$ws.Range("C20").ClearContents()

# Row 21
$ws.Range("A21").Value = 'Método:'
$ws.Range("B21").Value = 'O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas.'
$ws.Range("C21").Value = 'O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas.'

# Row 22
$ws.Range("A22").Value = 'Critério:'
$ws.Range("B22").Value = 'A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina.'
$ws.Range("C22").Value = 'A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina.'

# Row 23
$ws.Range("A23").Value = 'Norma de recuperação:'
$ws.Range("B23").Value = 'não há'
$ws.Range("C23").Value = 'não há'

# Row 24
$ws.Range("A24").Value = 'Bibliografia:'
$ws.Range("B24").Value = 'A ser definida em função do projeto'
$ws.Range("C24").Value = 'A ser definida em função do projeto'

# Row 25
$ws.Range("A25").Value = 'Requisitos:'
$ws.Range("B25").ClearContents()
$ws.Range("C25").ClearContents()

# Row 26
$ws.Range("A26").ClearContents()
$ws.Range("B26").Value = 'LOM3110 -  Projeto Integrado em Engenharia de Materiais III  (Requisito fraco)
'
$ws.Range("C26").Value = 'LOM3110 -  Projeto Integrado em Engenharia de Materiais III  (Requisito fraco)
'

# --- Row heights for rows 12-26 (restore default 15 or set the custom height) ---
$ws.Rows("12").RowHeight = 15
$ws.Rows("13").RowHeight = 15
$ws.Rows("14").RowHeight = 15
$ws.Rows("15").RowHeight = 15
$ws.Rows("16").RowHeight = 60
$ws.Rows("17").RowHeight = 15
$ws.Rows("18").RowHeight = 120
$ws.Rows("19").RowHeight = 120
$ws.Rows("20").RowHeight = 15
$ws.Rows("21").RowHeight = 60
$ws.Rows("22").RowHeight = 60
$ws.Rows("23").RowHeight = 60
$ws.Rows("24").RowHeight = 120
$ws.Rows("25").RowHeight = 15
$ws.Rows("26").RowHeight = 30

# --- Column layout: split the merged A:B column-width range so column A keeps its
#     original exact width while column B gets its own explicit width entry ---
$ws.Columns("B").ColumnWidth = $ws.Columns("B").ColumnWidth
